$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118 (Femacal de La Calera / Acelga weekly price entry),
# which pushes the existing rows 118-228 down to 119-229.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new weekly record.
$ws.Range("A118").Value = 3
$ws.Range("B118").Value = "Femacal de La Calera"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44512
$ws.Range("E118").Value = 5
$ws.Range("F118").Value = 100112009
$ws.Range("G118").Value = "Acelga"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 250
$ws.Range("K118").Value = 2000
$ws.Range("L118").Value = 2200
$ws.Range("M118").Value = 2104
$ws.Range("N118").Value = "$/docena de atados (6 kilos)"
$ws.Range("O118").Value = "Provincia de Quillota"
$ws.Range("P118").Value = 351
$ws.Range("Q118").Value = 6
$ws.Range("R118").Value = "Hortaliza"
